$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 15.02514266666667
$ws.Range("H2").Value = 45.075428
$ws.Range("I2").Value = 0.1401726531301337
$ws.Range("J2").Value = 0.1401726531301337
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 153.5290173333333
$ws.Range("N2").Value = 460.587052
$ws.Range("O2").Value = 0.3172206968818489
$ws.Range("P2").Value = 0.317220696881849
$ws.Range("Q2").Value = 2306.795388906473
$ws.Range("R2").Value = 20761.15850015826
$ws.Range("S2").Value = 0.04446566670971869
$ws.Range("T2").Value = 0.04446566670971869
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 15.02514266666667
$ws.Range("H3").Value = 45.075428
$ws.Range("I3").Value = 0.1401726531301337
$ws.Range("J3").Value = 0.1401726531301337
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 168.7997026666667
$ws.Range("N3").Value = 506.3991080000001
$ws.Range("O3").Value = 0.3487728915577651
$ws.Range("P3").Value = 0.3487728915577651
$ws.Range("Q3").Value = 2536.239614657581
$ws.Range("R3").Value = 22826.15653191823
$ws.Range("S3").Value = 0.04888842154952035
$ws.Range("T3").Value = 0.04888842154952034
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 15.02514266666667
$ws.Range("H4").Value = 45.075428
$ws.Range("I4").Value = 0.1401726531301337
$ws.Range("J4").Value = 0.1401726531301337
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 68.09032333333333
$ws.Range("N4").Value = 204.27097
$ws.Range("O4").Value = 0.1406878008722904
$ws.Range("P4").Value = 0.1406878008722904
$ws.Range("Q4").Value = 1023.066822302796
$ws.Range("R4").Value = 9207.601400725162
$ws.Range("S4").Value = 0.01972058231131288
$ws.Range("T4").Value = 0.01972058231131288
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 15.02514266666667
$ws.Range("H5").Value = 45.075428
$ws.Range("I5").Value = 0.1401726531301337
$ws.Range("J5").Value = 0.1401726531301337
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 93.562673
$ws.Range("N5").Value = 280.688019
$ws.Range("O5").Value = 0.1933186106880956
$ws.Range("P5").Value = 0.1933186106880956
$ws.Range("Q5").Value = 1405.792510099682
$ws.Range("R5").Value = 12652.13259089713
$ws.Range("S5").Value = 0.02709798255958179
$ws.Range("T5").Value = 0.02709798255958179
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 33.26311566666666
$ws.Range("H6").Value = 99.78934699999999
$ws.Range("I6").Value = 0.3103184627135109
$ws.Range("J6").Value = 0.3103184627135109
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 153.5290173333333
$ws.Range("N6").Value = 460.587052
$ws.Range("O6").Value = 0.3172206968818489
$ws.Range("P6").Value = 0.317220696881849
$ws.Range("Q6").Value = 5106.853461748337
$ws.Range("R6").Value = 45961.68115573504
$ws.Range("S6").Value = 0.09843943899728398
$ws.Range("T6").Value = 0.09843943899728398
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 33.26311566666666
$ws.Range("H7").Value = 99.78934699999999
$ws.Range("I7").Value = 0.3103184627135109
$ws.Range("J7").Value = 0.3103184627135109
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 168.7997026666667
$ws.Range("N7").Value = 506.3991080000001
$ws.Range("O7").Value = 0.3487728915577651
$ws.Range("P7").Value = 0.3487728915577651
$ws.Range("Q7").Value = 5614.804034300276
$ws.Range("R7").Value = 50533.23630870248
$ws.Range("S7").Value = 0.1082306675443517
$ws.Range("T7").Value = 0.1082306675443517
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 33.26311566666666
$ws.Range("H8").Value = 99.78934699999999
$ws.Range("I8").Value = 0.3103184627135109
$ws.Range("J8").Value = 0.3103184627135109
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 68.09032333333333
$ws.Range("N8").Value = 204.27097
$ws.Range("O8").Value = 0.1406878008722904
$ws.Range("P8").Value = 0.1406878008722904
$ws.Range("Q8").Value = 2264.896300817399
$ws.Range("R8").Value = 20384.06670735659
$ws.Range("S8").Value = 0.04365802208923369
$ws.Range("T8").Value = 0.04365802208923369
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 33.26311566666666
$ws.Range("H9").Value = 99.78934699999999
$ws.Range("I9").Value = 0.3103184627135109
$ws.Range("J9").Value = 0.3103184627135109
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 93.562673
$ws.Range("N9").Value = 280.688019
$ws.Range("O9").Value = 0.1933186106880956
$ws.Range("P9").Value = 0.1933186106880956
$ws.Range("Q9").Value = 3112.18601408151
$ws.Range("R9").Value = 28009.67412673359
$ws.Range("S9").Value = 0.05999033408264155
$ws.Range("T9").Value = 0.05999033408264153
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 50.12360066666667
$ws.Range("H10").Value = 150.370802
$ws.Range("I10").Value = 0.467613402797773
$ws.Range("J10").Value = 0.4676134027977729
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 153.5290173333333
$ws.Range("N10").Value = 460.587052
$ws.Range("O10").Value = 0.3172206968818489
$ws.Range("P10").Value = 0.317220696881849
$ws.Range("Q10").Value = 7695.427155561744
$ws.Range("R10").Value = 69258.84440005569
$ws.Range("S10").Value = 0.1483366495068023
$ws.Range("T10").Value = 0.1483366495068023
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 50.12360066666667
$ws.Range("H11").Value = 150.370802
$ws.Range("I11").Value = 0.467613402797773
$ws.Range("J11").Value = 0.4676134027977729
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 168.7997026666667
$ws.Range("N11").Value = 506.3991080000001
$ws.Range("O11").Value = 0.3487728915577651
$ws.Range("P11").Value = 0.3487728915577651
$ws.Range("Q11").Value = 8460.84888911607
$ws.Range("R11").Value = 76147.64000204463
$ws.Range("S11").Value = 0.1630908786249452
$ws.Range("T11").Value = 0.1630908786249451
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 50.12360066666667
$ws.Range("H12").Value = 150.370802
$ws.Range("I12").Value = 0.467613402797773
$ws.Range("J12").Value = 0.4676134027977729
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 68.09032333333333
$ws.Range("N12").Value = 204.27097
$ws.Range("O12").Value = 0.1406878008722904
$ws.Range("P12").Value = 0.1406878008722904
$ws.Range("Q12").Value = 3412.932176024216
$ws.Range("R12").Value = 30716.38958421794
$ws.Range("S12").Value = 0.0657875012980272
$ws.Range("T12").Value = 0.0657875012980272
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 50.12360066666667
$ws.Range("H13").Value = 150.370802
$ws.Range("I13").Value = 0.467613402797773
$ws.Range("J13").Value = 0.4676134027977729
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 93.562673
$ws.Range("N13").Value = 280.688019
$ws.Range("O13").Value = 0.1933186106880956
$ws.Range("P13").Value = 0.1933186106880956
$ws.Range("Q13").Value = 4689.698058757916
$ws.Range("R13").Value = 42207.28252882123
$ws.Range("S13").Value = 0.09039837336799833
$ws.Range("T13").Value = 0.0903983733679983
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 8.778397666666669
$ws.Range("H14").Value = 26.335193
$ws.Range("I14").Value = 0.08189548135858246
$ws.Range("J14").Value = 0.08189548135858243
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 153.5290173333333
$ws.Range("N14").Value = 460.587052
$ws.Range("O14").Value = 0.3172206968818489
$ws.Range("P14").Value = 0.317220696881849
$ws.Range("Q14").Value = 1347.73876752456
$ws.Range("R14").Value = 12129.64890772104
$ws.Range("S14").Value = 0.02597894166804399
$ws.Range("T14").Value = 0.02597894166804399
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 8.778397666666669
$ws.Range("H15").Value = 26.335193
$ws.Range("I15").Value = 0.08189548135858246
$ws.Range("J15").Value = 0.08189548135858243
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 168.7997026666667
$ws.Range("N15").Value = 506.3991080000001
$ws.Range("O15").Value = 0.3487728915577651
$ws.Range("P15").Value = 0.3487728915577651
$ws.Range("Q15").Value = 1481.790916023094
$ws.Range("R15").Value = 13336.11824420785
$ws.Range("S15").Value = 0.02856292383894785
$ws.Range("T15").Value = 0.02856292383894784
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 8.778397666666669
$ws.Range("H16").Value = 26.335193
$ws.Range("I16").Value = 0.08189548135858246
$ws.Range("J16").Value = 0.08189548135858243
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 68.09032333333333
$ws.Range("N16").Value = 204.27097
$ws.Range("O16").Value = 0.1406878008722904
$ws.Range("P16").Value = 0.1406878008722904
$ws.Range("Q16").Value = 597.7239354719123
$ws.Range("R16").Value = 5379.515419247211
$ws.Range("S16").Value = 0.01152169517371662
$ws.Range("T16").Value = 0.01152169517371662
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 8.778397666666669
$ws.Range("H17").Value = 26.335193
$ws.Range("I17").Value = 0.08189548135858246
$ws.Range("J17").Value = 0.08189548135858243
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 93.562673
$ws.Range("N17").Value = 280.688019
$ws.Range("O17").Value = 0.1933186106880956
$ws.Range("P17").Value = 0.1933186106880956
$ws.Range("Q17").Value = 821.3303503502965
$ws.Range("R17").Value = 7391.973153152668
$ws.Range("S17").Value = 0.015831920677874
$ws.Range("T17").Value = 0.015831920677874
